# Updated symbol list on Tue Dec 27 05:41:05 UTC 2022 with GitHub Actions
# Refresh the "Price" column (D) with the latest quotes and tweak a couple of
# "Volume(1h)" (E) labels that flipped their Best/Worst-in-24h badge.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> new price. These are text-formatted cells (not numbers), so
# force the Text number format before writing the value to keep them as
# strings instead of letting Excel auto-coerce them to doubles.
$priceUpdates = @{
    "D2"  = "243.34"
    "D3"  = "23.06"
    "D4"  = "5.403"
    "D5"  = "0.05986"
    "D6"  = "3.428"
    "D7"  = "6.490"
    "D8"  = "0.8080"
    "D9"  = "0.9235"
    "D10" = "0.1425"
    "D11" = "0.07417"
    "D12" = "0.03258"
    "D14" = "0.09361"
    "D15" = "3.854"
    "D16" = "0.001587"
    "D17" = "0.04700"
    "D18" = "0.0005890"
    "D19" = "0.005860"
    "D21" = "0.004880"
    "D22" = "0.00006800"
    "D23" = "3.573"
    "D27" = "0.0002340"
    "D40" = "0.03971"
    "D41" = "0.006350"
    "D42" = "0.004300"
    "D43" = "0.1077"
    "D44" = "0.008917"
    "D47" = "0.6500"
    "D48" = "0.002457"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Volume(1h) label tweaks (plain text, no coercion risk).
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
